# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap two pairs of province names (their row order changed after resort) ---
# Row 23 (Gipuzkoa/Guipuzcoa) <-> Row 24 (Valladolid)
$ws.Range("A23").Value = "Valladolid"
$ws.Range("A24").Value = "Gipuzkoa/Guipuzcoa"

# Row 28 (Caceres) <-> Row 29 (Segovia)
$ws.Range("A28").Value = "Segovia"
$ws.Range("A29").Value = "Caceres"

# --- Update "last updated" timestamp text ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 13:22"

# --- Update numeric data cells (Casos totales / Casos activos / Recuperados / Muertes) ---
# Row 20
$ws.Range("B20").Value = 1915
$ws.Range("C20").Value = 477
$ws.Range("D20").Value = 1213
$ws.Range("E20").Value = 225

# Row 22
$ws.Range("C22").Value = 294
$ws.Range("D22").Value = 1309

# Row 23
$ws.Range("B23").Value = 1686
$ws.Range("C23").Value = 604
$ws.Range("D23").Value = 920
$ws.Range("E23").Value = 162

# Row 24
$ws.Range("B24").Value = 1678
$ws.Range("C24").Value = 4151
$ws.Range("D24").Value = 4663
$ws.Range("E24").Value = 92

# Row 28
$ws.Range("B28").Value = 1480
$ws.Range("C28").Value = 415
$ws.Range("D28").Value = 943
$ws.Range("E28").Value = 122

# Row 29
$ws.Range("B29").Value = 1453
$ws.Range("C29").Value = 142
$ws.Range("D29").Value = 1083
$ws.Range("E29").Value = 228

# Row 30
$ws.Range("B30").Value = 1398
$ws.Range("C30").Value = 620
$ws.Range("D30").Value = 568
$ws.Range("E30").Value = 210

# Row 31
$ws.Range("C31").Value = 203
$ws.Range("D31").Value = 1038

# Row 34
$ws.Range("B34").Value = 1053
$ws.Range("C34").Value = 416
$ws.Range("D34").Value = 513
$ws.Range("E34").Value = 124

# Row 40
$ws.Range("B40").Value = 849
$ws.Range("C40").Value = 182
$ws.Range("E40").Value = 66

# Row 41
$ws.Range("B41").Value = 793
$ws.Range("C41").Value = 272
$ws.Range("D41").Value = 436
$ws.Range("E41").Value = 85

# Row 43
$ws.Range("B43").Value = 731
$ws.Range("C43").Value = 186
$ws.Range("D43").Value = 502
$ws.Range("E43").Value = 43

# Row 46
$ws.Range("B46").Value = 524
$ws.Range("C46").Value = 137
$ws.Range("D46").Value = 346

# Row 51
$ws.Range("B51").Value = 360
$ws.Range("C51").Value = 119
$ws.Range("D51").Value = 198
$ws.Range("E51").Value = 43
